$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Rename "Sheet1" -> "Optimizer"
$ws.Name = "Optimizer"

# --- Sheet content -----------------------------------------------------
# Text cells are written in this specific order so that new shared-string
# table entries land at the same indices the authored workbook used.
$ws.Range("B2").Value = "WMO1"
$ws.Range("C2").Value = "WMO2"
$ws.Range("D2").Value = "WMO3"
$ws.Range("A3").Value = "Volume"
$ws.Range("A4").Value = "Cost"
$ws.Range("A1").Value = "Shortage level:"
$ws.Range("A7").Value = "Upper Bounds"
$ws.Range("E2").Value = "Rationing"
$ws.Range("F2").Value = "Total Cost"

# Numbers / formulas
$ws.Range("B1").Value = 100

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 50
$ws.Range("D3").Value = 50
$ws.Range("E3").Formula = '=B1-SUM(B3:D3)'

$ws.Range("B4").Formula = '=B3*20'
$ws.Range("C4").Formula = '=C3*15'
$ws.Range("D4").Formula = '=D3*10'
$ws.Range("E4").Formula = '=E3*20'
$ws.Range("F4").Formula = '=SUM(B4:E4)'

$ws.Range("B7").Value = 50
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = 50

# --- Column widths (best-fit approximations) --------------------------
$ws.Columns.Item(1).ColumnWidth = 13.67
$ws.Columns.Item(4).ColumnWidth = 15
$ws.Columns.Item(5).ColumnWidth = 11.83
$ws.Columns.Item(9).ColumnWidth = 18.5

# --- Selection ---------------------------------------------------------
$ws.Range("A13").Select()

# --- Solver parameters (hidden, sheet-local defined names) -------------
function Add-SolverName($name, $formula) {
    $n = $ws.Names.Add($name, $formula)
    $n.Visible = $false
}

Add-SolverName "solver_adj" '=Optimizer!$B$3:$D$3'
Add-SolverName "solver_cvg" '=0.0001'
Add-SolverName "solver_drv" '=1'
Add-SolverName "solver_eng" '=1'
Add-SolverName "solver_est" '=1'
Add-SolverName "solver_itr" '=2147483647'
Add-SolverName "solver_lhs1" '=Optimizer!$B$4'
Add-SolverName "solver_lhs2" '=Optimizer!$C$3'
Add-SolverName "solver_lhs3" '=Optimizer!$D$3'
Add-SolverName "solver_lhs4" '=Optimizer!$E$3'
Add-SolverName "solver_lhs5" '=Optimizer!$F$4'
Add-SolverName "solver_mip" '=2147483647'
Add-SolverName "solver_mni" '=30'
Add-SolverName "solver_mrt" '=0.075'
Add-SolverName "solver_msl" '=2'
Add-SolverName "solver_neg" '=1'
Add-SolverName "solver_nod" '=2147483647'
Add-SolverName "solver_num" '=5'
Add-SolverName "solver_nwt" '=1'
Add-SolverName "solver_opt" '=Optimizer!$F$4'
Add-SolverName "solver_pre" '=0.000001'
Add-SolverName "solver_rbv" '=1'
Add-SolverName "solver_rel1" '=1'
Add-SolverName "solver_rel2" '=1'
Add-SolverName "solver_rel3" '=1'
Add-SolverName "solver_rel4" '=3'
Add-SolverName "solver_rel5" '=3'
Add-SolverName "solver_rhs1" '=Optimizer!$B$7'
Add-SolverName "solver_rhs2" '=Optimizer!$C$7'
Add-SolverName "solver_rhs3" '=Optimizer!$D$7'
Add-SolverName "solver_rhs4" '=0'
Add-SolverName "solver_rhs5" '=0'
Add-SolverName "solver_rlx" '=2'
Add-SolverName "solver_rsd" '=0'
Add-SolverName "solver_scl" '=1'
Add-SolverName "solver_sho" '=2'
Add-SolverName "solver_ssz" '=100'
Add-SolverName "solver_tim" '=2147483647'
Add-SolverName "solver_tol" '=0.01'
Add-SolverName "solver_typ" '=2'
Add-SolverName "solver_val" '=0'
Add-SolverName "solver_ver" '=3'

Write-Output "done"
